# Revert "Merging Katherine's working branch to main"
# Fix the run-level line breaks in the "Title 16" caption shape on slide 1:
#   Run 2: "An Empty Plot"          -> "An Empty Plot\n"
#   Run 3: "\nProject Description 1" -> "Project Description 1\n"
#   Run 4: "\nProject Description 2" -> "Project Description 2"
#
# Each run's text is rewritten in place via TextRange.Characters(start, length)
# so the existing run formatting (rPr) is preserved. Edits are applied from the
# last run back to the first so that the original (pre-edit) character offsets
# stay valid for every Characters() call.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Run 4 (originally chars 47-68): "\nProject Description 2" -> "Project Description 2"
$run4 = $tr.Characters(47, 22)
$run4.Text = "Project Description 2"

# Run 3 (originally chars 25-46): "\nProject Description 1" -> "Project Description 1\n"
$run3 = $tr.Characters(25, 22)
$run3.Text = "Project Description 1`n"

# Run 2 (originally chars 12-24): "An Empty Plot" -> "An Empty Plot\n"
$run2 = $tr.Characters(12, 13)
$run2.Text = "An Empty Plot`n"
